$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new Price text would otherwise be auto-parsed as a number
# (e.g. "1.0000" -> 1). Mark them as Text first so the literal string is kept,
# matching the inline-string cells produced by the source data refresh.
$textCells = @("D4", "D5", "D7", "D9", "D12", "D13", "D14", "D15", "D17", "D18", "D19", "D20", "D22", "D23", "D24", "D25", "D27", "D28", "D29", "D30", "D31", "D32", "D33", "D35", "D36", "D37", "D38", "D39", "D40", "D41", "D43", "D44", "D47", "D48", "D49", "D50", "D51")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = '28.704.24'
$ws.Range("E2").Value = '  +6.78%  '
$ws.Range("D3").Value = '1.808.87'
$ws.Range("E3").Value = '  +4.54%  '
$ws.Range("D4").Value = '1.0000'
$ws.Range("E4").Value = '  +0.27%  '
$ws.Range("D5").Value = '251.94'
$ws.Range("E5").Value = '  +4.09%  '
$ws.Range("E6").Value = '  +0.26%  '
$ws.Range("D7").Value = '0.4966'
$ws.Range("E7").Value = '  +1.03%  '
$ws.Range("E8").Value = '  +7.18%  '
$ws.Range("D9").Value = '0.06383'
$ws.Range("E9").Value = '  +2.56%  '
$ws.Range("D10").Value = '1.803.65'
$ws.Range("E10").Value = '  +4.17%  '
$ws.Range("E11").Value = '  +4.50%  '
$ws.Range("D12").Value = '0.07115'
$ws.Range("E12").Value = '  +3.14%  '
$ws.Range("D13").Value = '0.6478'
$ws.Range("E13").Value = '  +6.11%  '
$ws.Range("D14").Value = '4.707'
$ws.Range("E14").Value = '  +4.52%  '
$ws.Range("D15").Value = '81.90'
$ws.Range("E15").Value = '  +5.88%  '
$ws.Range("D16").Value = '28.684.72'
$ws.Range("E16").Value = '  +6.75%  '
$ws.Range("D17").Value = '1.0000'
$ws.Range("E17").Value = '  +0.19%  '
$ws.Range("D18").Value = '0.000007354'
$ws.Range("E18").Value = '  +2.32%  '
$ws.Range("D19").Value = '0.9998'
$ws.Range("E19").Value = '  +0.25%  '
$ws.Range("D20").Value = '12.26'
$ws.Range("E20").Value = '  +6.95%  '
$ws.Range("D21").Value = '2.036.92'
$ws.Range("E21").Value = '  +3.97%  '
$ws.Range("D22").Value = '4.619'
$ws.Range("E22").Value = '  +4.08%  '
$ws.Range("D23").Value = '8.892'
$ws.Range("E23").Value = '  +3.84%  '
$ws.Range("D24").Value = '5.313'
$ws.Range("E24").Value = '  +3.65%  '
$ws.Range("D25").Value = '142.85'
$ws.Range("E25").Value = '  +2.90%  '
$ws.Range("E26").Value = '  +4.48%  '
$ws.Range("D27").Value = '1.881'
$ws.Range("E27").Value = '  +4.97%  '
$ws.Range("D28").Value = '111.25'
$ws.Range("E28").Value = '  +4.69%  '
$ws.Range("D29").Value = '1.387'
$ws.Range("E29").Value = '  +0.40%  '
$ws.Range("D30").Value = '4.185'
$ws.Range("E30").Value = '  +5.98%  '
$ws.Range("D31").Value = '0.08361'
$ws.Range("E31").Value = '  +4.58%  '
$ws.Range("D32").Value = '3.845'
$ws.Range("E32").Value = '  +4.47%  '
$ws.Range("D33").Value = '0.04964'
$ws.Range("E33").Value = '  +9.44%  '
$ws.Range("E34").Value = '  +8.07%  '
$ws.Range("D35").Value = '0.6725'
$ws.Range("E35").Value = '  +7.67%  '
$ws.Range("D36").Value = '2.666'
$ws.Range("E36").Value = '  +2.13%  '
$ws.Range("D37").Value = '0.9601'
$ws.Range("E37").Value = '  +2.90%  '
$ws.Range("D38").Value = '2.643'
$ws.Range("E38").Value = '  +8.33%  '
$ws.Range("D39").Value = '2.147'
$ws.Range("E39").Value = '  +4.17%  '
$ws.Range("D40").Value = '0.01598'
$ws.Range("E40").Value = '  +6.24%  '
$ws.Range("D41").Value = '5.938'
$ws.Range("E41").Value = '  +5.01%  '
$ws.Range("E42").Value = '  -0.06%  '
$ws.Range("D43").Value = '101.04'
$ws.Range("E43").Value = '  +1.25%  '
$ws.Range("D44").Value = '0.4120'
$ws.Range("E44").Value = '  +6.40%  '
$ws.Range("E45").Value = '  +4.14%  '
$ws.Range("E46").Value = '  +5.52%  '
$ws.Range("D47").Value = '0.05492'
$ws.Range("E47").Value = '  +1.91%  '
$ws.Range("D48").Value = '8.159'
$ws.Range("E48").Value = '  +2.56%  '
$ws.Range("D49").Value = '31.38'
$ws.Range("E49").Value = '  +3.74%  '
$ws.Range("D50").Value = '1.304'
$ws.Range("E50").Value = '  +4.98%  '
$ws.Range("D51").Value = '0.3606'
$ws.Range("E51").Value = '  +6.46%  '
